$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.222587333333333
$ws.Range("H2").Value = 6.667762
$ws.Range("I2").Value = 0.1134117015526119
$ws.Range("J2").Value = 0.1134117015526119
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 35.31114333333333
$ws.Range("N2").Value = 105.93343
$ws.Range("O2").Value = 0.6187867769880316
$ws.Range("P2").Value = 0.6187867769880316
$ws.Range("Q2").Value = 78.48209989818444
$ws.Range("R2").Value = 706.33889908366
$ws.Range("S2").Value = 0.0701776612764693
$ws.Range("T2").Value = 0.07017766127646928

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.222587333333333
$ws.Range("H3").Value = 6.667762
$ws.Range("I3").Value = 0.1134117015526119
$ws.Range("J3").Value = 0.1134117015526119
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.551362
$ws.Range("N3").Value = 19.654086
$ws.Range("O3").Value = 0.1148050103785518
$ws.Range("P3").Value = 0.1148050103785518
$ws.Range("Q3").Value = 14.56097419728133
$ws.Range("R3").Value = 131.048767775532
$ws.Range("S3").Value = 0.01302023157379684
$ws.Range("T3").Value = 0.01302023157379683

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.222587333333333
$ws.Range("H4").Value = 6.667762
$ws.Range("I4").Value = 0.1134117015526119
$ws.Range("J4").Value = 0.1134117015526119
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.951915
$ws.Range("N4").Value = 32.855745
$ws.Range("O4").Value = 0.191919590955288
$ws.Range("P4").Value = 0.191919590955288
$ws.Range("Q4").Value = 24.34158755474333
$ws.Range("R4").Value = 219.07428799269
$ws.Range("S4").Value = 0.02176592737152049
$ws.Range("T4").Value = 0.02176592737152048

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.222587333333333
$ws.Range("H5").Value = 6.667762
$ws.Range("I5").Value = 0.1134117015526119
$ws.Range("J5").Value = 0.1134117015526119
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.250702333333333
$ws.Range("N5").Value = 12.752107
$ws.Range("O5").Value = 0.07448862167812857
$ws.Range("P5").Value = 0.07448862167812857
$ws.Range("Q5").Value = 9.447557163837109
$ws.Range("R5").Value = 85.02801447453399
$ws.Range("S5").Value = 0.008447881330825338
$ws.Range("T5").Value = 0.008447881330825337

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.89424733333333
$ws.Range("H6").Value = 35.682742
$ws.Range("I6").Value = 0.6069263549423107
$ws.Range("J6").Value = 0.6069263549423106
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.31114333333333
$ws.Range("N6").Value = 105.93343
$ws.Range("O6").Value = 0.6187867769880316
$ws.Range("P6").Value = 0.6187867769880316
$ws.Range("Q6").Value = 419.9994724294511
$ws.Range("R6").Value = 3779.99525186506
$ws.Range("S6").Value = 0.3755580030438466
$ws.Range("T6").Value = 0.3755580030438465

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.89424733333333
$ws.Range("H7").Value = 35.682742
$ws.Range("I7").Value = 0.6069263549423107
$ws.Range("J7").Value = 0.6069263549423106
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.551362
$ws.Range("N7").Value = 19.654086
$ws.Range("O7").Value = 0.1148050103785518
$ws.Range("P7").Value = 0.1148050103785518
$ws.Range("Q7").Value = 77.92351999820133
$ws.Range("R7").Value = 701.3116799838119
$ws.Range("S7").Value = 0.06967818647816862
$ws.Range("T7").Value = 0.0696781864781686

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.89424733333333
$ws.Range("H8").Value = 35.682742
$ws.Range("I8").Value = 0.6069263549423107
$ws.Range("J8").Value = 0.6069263549423106
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.951915
$ws.Range("N8").Value = 32.855745
$ws.Range("O8").Value = 0.191919590955288
$ws.Range("P8").Value = 0.191919590955288
$ws.Range("Q8").Value = 130.2647857836433
$ws.Range("R8").Value = 1172.38307205279
$ws.Range("S8").Value = 0.1164810577805122
$ws.Range("T8").Value = 0.1164810577805122

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.89424733333333
$ws.Range("H9").Value = 35.682742
$ws.Range("I9").Value = 0.6069263549423107
$ws.Range("J9").Value = 0.6069263549423106
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.250702333333333
$ws.Range("N9").Value = 12.752107
$ws.Range("O9").Value = 0.07448862167812857
$ws.Range("P9").Value = 0.07448862167812857
$ws.Range("Q9").Value = 50.55890489304377
$ws.Range("R9").Value = 455.0301440373939
$ws.Range("S9").Value = 0.04520910763978336
$ws.Range("T9").Value = 0.04520910763978336

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.762
$ws.Range("H10").Value = 2.286
$ws.Range("I10").Value = 0.0388824840702579
$ws.Range("J10").Value = 0.03888248407025789
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 35.31114333333333
$ws.Range("N10").Value = 105.93343
$ws.Range("O10").Value = 0.6187867769880316
$ws.Range("P10").Value = 0.6187867769880316
$ws.Range("Q10").Value = 26.90709122
$ws.Range("R10").Value = 242.16382098
$ws.Range("S10").Value = 0.02405996699912337
$ws.Range("T10").Value = 0.02405996699912336

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.762
$ws.Range("H11").Value = 2.286
$ws.Range("I11").Value = 0.0388824840702579
$ws.Range("J11").Value = 0.03888248407025789
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.551362
$ws.Range("N11").Value = 19.654086
$ws.Range("O11").Value = 0.1148050103785518
$ws.Range("P11").Value = 0.1148050103785518
$ws.Range("Q11").Value = 4.992137844
$ws.Range("R11").Value = 44.929240596
$ws.Range("S11").Value = 0.004463903987229834
$ws.Range("T11").Value = 0.004463903987229833

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.762
$ws.Range("H12").Value = 2.286
$ws.Range("I12").Value = 0.0388824840702579
$ws.Range("J12").Value = 0.03888248407025789
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.951915
$ws.Range("N12").Value = 32.855745
$ws.Range("O12").Value = 0.191919590955288
$ws.Range("P12").Value = 0.191919590955288
$ws.Range("Q12").Value = 8.34535923
$ws.Range("R12").Value = 75.10823307
$ws.Range("S12").Value = 0.007462310438089397
$ws.Range("T12").Value = 0.007462310438089396

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.762
$ws.Range("H13").Value = 2.286
$ws.Range("I13").Value = 0.0388824840702579
$ws.Range("J13").Value = 0.03888248407025789
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.250702333333333
$ws.Range("N13").Value = 12.752107
$ws.Range("O13").Value = 0.07448862167812857
$ws.Range("P13").Value = 0.07448862167812857
$ws.Range("Q13").Value = 3.239035178
$ws.Range("R13").Value = 29.151316602
$ws.Range("S13").Value = 0.002896302645815301
$ws.Range("T13").Value = 0.002896302645815301

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.718678666666666
$ws.Range("H14").Value = 14.156036
$ws.Range("I14").Value = 0.2407794594348195
$ws.Range("J14").Value = 0.2407794594348195
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 35.31114333333333
$ws.Range("N14").Value = 105.93343
$ws.Range("O14").Value = 0.6187867769880316
$ws.Range("P14").Value = 0.6187867769880316
$ws.Range("Q14").Value = 166.6219387426089
$ws.Range("R14").Value = 1499.59744868348
$ws.Range("S14").Value = 0.1489911456685925
$ws.Range("T14").Value = 0.1489911456685925

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.718678666666666
$ws.Range("H15").Value = 14.156036
$ws.Range("I15").Value = 0.2407794594348195
$ws.Range("J15").Value = 0.2407794594348195
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 6.551362
$ws.Range("N15").Value = 19.654086
$ws.Range("O15").Value = 0.1148050103785518
$ws.Range("P15").Value = 0.1148050103785518
$ws.Range("Q15").Value = 30.91377210701067
$ws.Range("R15").Value = 278.223948963096
$ws.Range("S15").Value = 0.02764268833935655
$ws.Range("T15").Value = 0.02764268833935655

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.718678666666666
$ws.Range("H16").Value = 14.156036
$ws.Range("I16").Value = 0.2407794594348195
$ws.Range("J16").Value = 0.2407794594348195
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.951915
$ws.Range("N16").Value = 32.855745
$ws.Range("O16").Value = 0.191919590955288
$ws.Range("P16").Value = 0.191919590955288
$ws.Range("Q16").Value = 51.67856766964666
$ws.Range("R16").Value = 465.10710902682
$ws.Range("S16").Value = 0.04621029536516591
$ws.Range("T16").Value = 0.04621029536516591

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.718678666666666
$ws.Range("H17").Value = 14.156036
$ws.Range("I17").Value = 0.2407794594348195
$ws.Range("J17").Value = 0.2407794594348195
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.250702333333333
$ws.Range("N17").Value = 12.752107
$ws.Range("O17").Value = 0.07448862167812857
$ws.Range("P17").Value = 0.07448862167812857
$ws.Range("Q17").Value = 20.05769841865022
$ws.Range("R17").Value = 180.519285767852
$ws.Range("S17").Value = 0.01793533006170457
$ws.Range("T17").Value = 0.01793533006170457
